$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 329 (everything from row 329 down shifts down by 4 rows,
# which naturally pushes the former last weekly group - rows 450:453 - out to rows 454:457).
$ws.Rows("329:332").Insert()

# Fill in the new rows with the new weekly price-report group (fecha = 44755),
# using the same constant columns (A,B,C,E,F,G,H,N,O,Q,R) as the rest of this
# "Comercializadora del Agro de Limari - Pepino dulce" block.
$rows = 329..332
$calidad = @("Especial", "Primera", "Segunda", "Tercera")
$volumen = @(500, 700, 500, 400)
$pmin    = @(12000, 10000, 7000, 4000)
$pmax    = @(13000, 11000, 8000, 5000)
$pprom   = @(12500, 10500, 7500, 4500)
$pkg     = @(694, 583, 417, 250)

for ($i = 0; $i -lt 4; $i++) {
    $r = $rows[$i]

    $ws.Cells.Item($r, 1).Value = 2
    $ws.Cells.Item($r, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44755
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = 100112043
    $ws.Cells.Item($r, 7).Value = "Pepino dulce"
    $ws.Cells.Item($r, 8).Value = "Cultivar IV Región"
    $ws.Cells.Item($r, 9).Value = $calidad[$i]
    $ws.Cells.Item($r, 10).Value = $volumen[$i]
    $ws.Cells.Item($r, 11).Value = $pmin[$i]
    $ws.Cells.Item($r, 12).Value = $pmax[$i]
    $ws.Cells.Item($r, 13).Value = $pprom[$i]
    $ws.Cells.Item($r, 14).Value = '$/bandeja 18 kilos'
    $ws.Cells.Item($r, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 16).Value = $pkg[$i]
    $ws.Cells.Item($r, 17).Value = 18
    $ws.Cells.Item($r, 18).Value = "Hortaliza"

    # Match the date-style used by the rest of the "Fecha" column (D).
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(333, 4).NumberFormat
}
